$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.093.52"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "3.315.93"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.48"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.03"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "3.893.73"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.133"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.68"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "66.198.09"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "3.299.40"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "434.14"
$ws.Range("E18").Value = "  -3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.19"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "3.471.89"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.193"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.21"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.67"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").Value = "2.843.89"
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.789"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.16"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.45"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0666"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.08"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "326.43"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0272"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.14"
$ws.Range("E51").Value = "  -1.27%  "
